$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 change
$ws.Range("G2").Value = 1.71

# Row 4 changes
$ws.Range("J4").Value = 4.8
$ws.Range("K4").Value = 5.2
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 3.2
$ws.Range("AE4").Value = 19
